# Update "想去人数" (want-to-go) counts on three worksheets to match
# the newly generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 32
$ws1.Range("F8").Value = 8076
$ws1.Range("F11").Value = 1100
$ws1.Range("F12").Value = 791
$ws1.Range("F13").Value = 38
$ws1.Range("F16").Value = 63
$ws1.Range("F19").Value = 859

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 23

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 32
$ws4.Range("F9").Value = 8076
$ws4.Range("F12").Value = 1100
$ws4.Range("F13").Value = 791
$ws4.Range("F14").Value = 38
$ws4.Range("F17").Value = 63
$ws4.Range("F20").Value = 859
$ws4.Range("F21").Value = 23
